$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 new rows after row 16 to extend the sheet to 19 rows,
# then copy formatting from row 16 (A column style with border/bold/center) to new A17:A19
$ws.Range("A16").Copy() | Out-Null
$ws.Range("A17:A19").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

# --- Update rows 10-19 (data for averaging schemes, now incl. spiral quadrature schemes) ---
$ws.Cells.Item(10, 1).Value2 = 8
$ws.Cells.Item(10, 2).Value2 = "Gaussian-Quadrature"
$ws.Cells.Item(10, 3).Value2 = 0.985526333062672
$ws.Cells.Item(10, 4).Value2 = 1.045570974701326
$ws.Cells.Item(10, 5).Value2 = 0.9877485424670935
$ws.Cells.Item(10, 6).Value2 = 0.985526333062672
$ws.Cells.Item(10, 7).Value2 = 1.026035405532553
$ws.Cells.Item(10, 8).Value2 = 0.9671672344821337
$ws.Cells.Item(10, 9).Value2 = 0.9877192943386085
$ws.Cells.Item(10, 10).Value2 = 1.045570974701326
$ws.Cells.Item(10, 11).Value2 = 1.01665975858421
$ws.Cells.Item(10, 12).Value2 = 1.001093045823441
$ws.Cells.Item(10, 13).Value2 = 0.9999612974307311

$ws.Cells.Item(11, 1).Value2 = 9
$ws.Cells.Item(11, 2).Value2 = "Spiral-90deg-10rot-5space"
$ws.Cells.Item(11, 3).Value2 = 1.002606832088758
$ws.Cells.Item(11, 4).Value2 = 0.974579858789121
$ws.Cells.Item(11, 5).Value2 = 1.005509305914233
$ws.Cells.Item(11, 6).Value2 = 1.002606832088758
$ws.Cells.Item(11, 7).Value2 = 0.9831367361760647
$ws.Cells.Item(11, 8).Value2 = 1.01733271551605
$ws.Cells.Item(11, 9).Value2 = 1.0058271505838
$ws.Cells.Item(11, 10).Value2 = 0.974579858789121
$ws.Cells.Item(11, 11).Value2 = 0.9900445823516771
$ws.Cells.Item(11, 12).Value2 = 0.9963257072202176
$ws.Cells.Item(11, 13).Value2 = 0.9981654331780044

$ws.Cells.Item(12, 1).Value2 = 10
$ws.Cells.Item(12, 2).Value2 = "Spiral-90deg-15rot-5space"
$ws.Cells.Item(12, 3).Value2 = 1.002579146304509
$ws.Cells.Item(12, 4).Value2 = 0.9746684457516152
$ws.Cells.Item(12, 5).Value2 = 1.005490925928601
$ws.Cells.Item(12, 6).Value2 = 1.002579146304509
$ws.Cells.Item(12, 7).Value2 = 0.9831551044872199
$ws.Cells.Item(12, 8).Value2 = 1.017314756692279
$ws.Cells.Item(12, 9).Value2 = 1.005824412574981
$ws.Cells.Item(12, 10).Value2 = 0.9746684457516152
$ws.Cells.Item(12, 11).Value2 = 0.9900796858401082
$ws.Cells.Item(12, 12).Value2 = 0.9963294160723084
$ws.Cells.Item(12, 13).Value2 = 0.9981721319565343

$ws.Cells.Item(13, 1).Value2 = 11
$ws.Cells.Item(13, 2).Value2 = "Spiral-90deg-10rot-3space"
$ws.Cells.Item(13, 3).Value2 = 1.002594410215886
$ws.Cells.Item(13, 4).Value2 = 0.9746178856888805
$ws.Cells.Item(13, 5).Value2 = 1.005498860437817
$ws.Cells.Item(13, 6).Value2 = 1.002594410215886
$ws.Cells.Item(13, 7).Value2 = 0.9831233690071572
$ws.Cells.Item(13, 8).Value2 = 1.017292713152137
$ws.Cells.Item(13, 9).Value2 = 1.00583564428801
$ws.Cells.Item(13, 10).Value2 = 0.9746178856888805
$ws.Cells.Item(13, 11).Value2 = 0.9900583730633487
$ws.Cells.Item(13, 12).Value2 = 0.9963263916396173
$ws.Cells.Item(13, 13).Value2 = 0.9981604804649812

$ws.Cells.Item(14, 1).Value2 = 12
$ws.Cells.Item(14, 2).Value2 = "NoRotation-tilt60deg"
$ws.Cells.Item(14, 3).Value2 = 0.9706200000000003
$ws.Cells.Item(14, 4).Value2 = 1.102935999999999
$ws.Cells.Item(14, 5).Value2 = 0.9733280000000002
$ws.Cells.Item(14, 6).Value2 = 0.9706200000000003
$ws.Cells.Item(14, 7).Value2 = 1.057511999999998
$ws.Cells.Item(14, 8).Value2 = 0.9422840000000009
$ws.Cells.Item(14, 9).Value2 = 0.9733279999999995
$ws.Cells.Item(14, 10).Value2 = 1.102935999999999
$ws.Cells.Item(14, 11).Value2 = 1.038132
$ws.Cells.Item(14, 12).Value2 = 1.004376
$ws.Cells.Item(14, 13).Value2 = 1.003334666666666

$ws.Cells.Item(15, 1).Value2 = 13
$ws.Cells.Item(15, 2).Value2 = "Rotation-NoTilt"
$ws.Cells.Item(15, 3).Value2 = 0.95
$ws.Cells.Item(15, 4).Value2 = 1.19375
$ws.Cells.Item(15, 5).Value2 = 0.95
$ws.Cells.Item(15, 6).Value2 = 0.95
$ws.Cells.Item(15, 7).Value2 = 1.11
$ws.Cells.Item(15, 8).Value2 = 0.89
$ws.Cells.Item(15, 9).Value2 = 0.95
$ws.Cells.Item(15, 10).Value2 = 1.19375
$ws.Cells.Item(15, 11).Value2 = 1.071875
$ws.Cells.Item(15, 12).Value2 = 1.0109375
$ws.Cells.Item(15, 13).Value2 = 1.007291666666667

$ws.Cells.Item(16, 1).Value2 = 14
$ws.Cells.Item(16, 2).Value2 = "Rotation-60detTilt"
$ws.Cells.Item(16, 3).Value2 = 0.9703958495232021
$ws.Cells.Item(16, 4).Value2 = 1.112147764531201
$ws.Cells.Item(16, 5).Value2 = 0.9707931623424055
$ws.Cells.Item(16, 6).Value2 = 0.9703958495232021
$ws.Cells.Item(16, 7).Value2 = 1.062469835775997
$ws.Cells.Item(16, 8).Value2 = 0.9353466490880022
$ws.Cells.Item(16, 9).Value2 = 0.9714232160256024
$ws.Cells.Item(16, 10).Value2 = 1.112147764531201
$ws.Cells.Item(16, 11).Value2 = 1.041470463436803
$ws.Cells.Item(16, 12).Value2 = 1.005933156480003
$ws.Cells.Item(16, 13).Value2 = 1.003762746214402

$ws.Cells.Item(17, 1).Value2 = 15
$ws.Cells.Item(17, 2).Value2 = "HexGrid-90degTilt5degRes"
$ws.Cells.Item(17, 3).Value2 = 0.9993414994895774
$ws.Cells.Item(17, 4).Value2 = 0.9986410695512723
$ws.Cells.Item(17, 5).Value2 = 0.999060179056716
$ws.Cells.Item(17, 6).Value2 = 0.9993414994895774
$ws.Cells.Item(17, 7).Value2 = 0.9989486674920069
$ws.Cells.Item(17, 8).Value2 = 0.9984177122683611
$ws.Cells.Item(17, 9).Value2 = 0.9990244158766295
$ws.Cells.Item(17, 10).Value2 = 0.9986410695512723
$ws.Cells.Item(17, 11).Value2 = 0.9988506243039941
$ws.Cells.Item(17, 12).Value2 = 0.9990960618967858
$ws.Cells.Item(17, 13).Value2 = 0.9989055906224271

$ws.Cells.Item(18, 1).Value2 = 16
$ws.Cells.Item(18, 2).Value2 = "HexGrid-90degTilt22p5degRes"
$ws.Cells.Item(18, 3).Value2 = 1.000058882978371
$ws.Cells.Item(18, 4).Value2 = 0.9919224248575809
$ws.Cells.Item(18, 5).Value2 = 1.001150001574911
$ws.Cells.Item(18, 6).Value2 = 1.000058882978371
$ws.Cells.Item(18, 7).Value2 = 0.9948337630093237
$ws.Cells.Item(18, 8).Value2 = 1.002875597164246
$ws.Cells.Item(18, 9).Value2 = 0.999967641797497
$ws.Cells.Item(18, 10).Value2 = 0.9919224248575809
$ws.Cells.Item(18, 11).Value2 = 0.9965362132162461
$ws.Cells.Item(18, 12).Value2 = 0.9982975480973084
$ws.Cells.Item(18, 13).Value2 = 0.9984680518969883

$ws.Cells.Item(19, 1).Value2 = 17
$ws.Cells.Item(19, 2).Value2 = "HexGrid-60degTilt5degRes"
$ws.Cells.Item(19, 3).Value2 = 1.003639802234912
$ws.Cells.Item(19, 4).Value2 = 0.9805524437189224
$ws.Cells.Item(19, 5).Value2 = 1.003976604818678
$ws.Cells.Item(19, 6).Value2 = 1.003639802234912
$ws.Cells.Item(19, 7).Value2 = 0.9878696545158564
$ws.Cells.Item(19, 8).Value2 = 1.012062097201732
$ws.Cells.Item(19, 9).Value2 = 1.004813196395894
$ws.Cells.Item(19, 10).Value2 = 0.9805524437189224
$ws.Cells.Item(19, 11).Value2 = 0.9922645242688002
$ws.Cells.Item(19, 12).Value2 = 0.9979521632518563
$ws.Cells.Item(19, 13).Value2 = 0.9988189664809992

